# Added and Updated Czech Test Data
#
# The Belgium sheet's full used range gets selected (no longer the active
# tab), and a new "Czech" sheet is created (as a copy of Belgium, since it
# shares the same "Market"/NGC values) and placed after it, becoming the
# active tab with cell C7 selected.

$wb = $excel.ActiveWorkbook

$belgium = $wb.Worksheets.Item("Belgium")

# Belgium is no longer the active tab; its full data range ends up selected.
$belgium.Range("A1:D16").Select()

# Create the new "Czech" sheet at the end of the workbook (after Belgium),
# by copying Belgium so the formatting/values/styles line up, then rename.
$belgium.Copy($null, $belgium)
$czech = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Name = "Czech"

# Reset row 16 back to the sheet's default (non-custom) height.
$czech.Rows.Item(16).AutoFit()

# Match the recorded column widths for the new sheet as closely as possible.
$czech.Columns.Item(2).ColumnWidth = 14.25
$czech.Columns.Item(3).ColumnWidth = 12.75
$czech.Columns.Item(4).ColumnWidth = 14.6

# Czech becomes the active tab with C7 selected.
$czech.Range("C7").Select()
